# "cambios en los datos" - update the IBAN list: keep only the header and
# the "ES32 5829 354  1000" entry (which moves up into row 2), removing
# the other four IBAN rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ES32 5829 354  1000"
$ws.Range("A3:A6").ClearContents() | Out-Null

$ws.Range("A2").Select() | Out-Null
